$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds numeric-looking text (e.g. "249.77"); Excel
# auto-converts such strings to numbers on assignment, so for those cells we
# force Text format, assign, then restore the default "Normal" style so the
# saved cell matches the original (unstyled) text cell.

$ws.Range("D2").Value = "37.142.95"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "2.053.00"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.672"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.34"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.54%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0794"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "

$ws.Range("E11").Value = "  +1.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.89%  "

$ws.Range("D13").Value = "2.349.92"
$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.833"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.73%  "

$ws.Range("D16").Value = "2.051.34"
$ws.Range("E16").Value = "  -0.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +29.91%  "

$ws.Range("D18").Value = "37.109.83"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "76.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.11%  "

$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  -4.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.24%  "

$ws.Range("E25").Value = "  +11.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.127"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "

$ws.Range("E30").Value = "  +7.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0631"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0885"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.89%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("E40").Value = "  +11.59%  "

$ws.Range("E41").Value = "  +20.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.27%  "

$ws.Range("D48").Value = "1.293.68"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").Value = "2.241.08"
$ws.Range("E51").Value = "  -0.43%  "
